$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 5, pushing the existing weekly records (old rows 5-16)
# down to rows 6-17 and picking up the date-format styling from the row above.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value = "Metropolitana"
$ws.Cells.Item(5, 4).Value = 44980
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100101
$ws.Cells.Item(5, 8).Value = "Berries"
$ws.Cells.Item(5, 9).Value = 100101008
$ws.Cells.Item(5, 10).Value = "Mora"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 250
$ws.Cells.Item(5, 14).Value = 4000
$ws.Cells.Item(5, 15).Value = 4000
$ws.Cells.Item(5, 16).Value = 4000
$ws.Cells.Item(5, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(5, 19).Value = 2000
$ws.Cells.Item(5, 20).Value = 2
